$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 73, shifting existing data (rows 73-205) down to (74-206)
$ws.Rows.Item(73).Insert()

# Populate the new row 73 with values. Most columns repeat the constant pattern
# used throughout the sheet; D, J, K, L, M, P hold the new record's data.
$ws.Range("A73").Value = 8
$ws.Range("B73").Value = "Terminal La Palmera de La Serena"
$ws.Range("C73").Value = "Coquimbo"
$ws.Range("D73").Value = 45259
$ws.Range("D73").Style = $ws.Range("D74").Style
$ws.Range("D73").NumberFormat = $ws.Range("D74").NumberFormat
$ws.Range("E73").Value = 4
$ws.Range("F73").Value = 100114007
$ws.Range("G73").Value = "Jengibre"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 360
$ws.Range("K73").Value = 23000
$ws.Range("L73").Value = 24000
$ws.Range("M73").Value = 23500
$ws.Range("N73").Value = "$/caja 13 kilos"
$ws.Range("O73").Value = "Perú"
$ws.Range("P73").Value = 1808
$ws.Range("Q73").Value = 13
$ws.Range("R73").Value = "Hortaliza"

Write-Host "done"
